# "final update for today"
# Append the new "Tenure" (ACS table B25007) rows to the Sheet1 data
# dictionary and move the on-screen selection down to the newly added
# rows, mirroring the author's last edit to Data Guide.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A=ID, B=Data Point, C=Column Name, D=Sub Category, E=ACS Variable.
# The sheet already has rows 2..114 (ID 1..113); we extend it with ID
# 114..122 (rows 115..123).
$newRows = @(
    @(114, "Tenure Total Households Series",  "tenure_total_agehh_series", "Tenure", "B25007_001E"),
    @(115, "Tenure Total Owner Occupied",      "tenure_allowneroccupied",   "Tenure", "B25007_002E"),
    @(116, "Tenure Owner Occupied 65 to 74",   "tenure_owner_65to74",       "Tenure", "B25007_009E"),
    @(117, "Tenure Owner Occupied 75 to 84",   "tenure_owner_75to84",       "Tenure", "B25007_010E"),
    @(118, "Tenure Owner Occupied 85+",        "tenure_owner_85+",          "Tenure", "B25007_011E"),
    @(119, "Tenure Total Renter Occupied",     "tenure_allrenteroccupied",  "Tenure", "B25007_012E"),
    @(120, "Tenure Renter Occupied 65 to 74",  "tenure_renter_65to74",      "Tenure", "B25007_019E"),
    @(121, "Tenure Renter Occupied 75 to 84",  "tenure_renter_75to84",      "Tenure", "B25007_020E"),
    @(122, "Tenure Renter Occupied 85+",       "tenure_renter_85+",         "Tenure", "B25007_021E")
)

$startRow = 115
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

$lastRow = $startRow + $newRows.Count - 1

# Match the author's final view/selection state: scrolled down to the new
# block and C115:C123 selected (was C40:C60, topLeftCell A82, before the
# edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
$selRange = $ws.Range($ws.Cells.Item($startRow, 3), $ws.Cells.Item($lastRow, 3))
$selRange.Select()
